{"js": "// Replace the date line and the 25 division problems with their updated\n// values, as described by the diff. Every <w:t> run in the document is\n// mapped 1:1 (old text -> new text), so a simple search-and-replace per\n// pair is sufficient and keeps each run's original formatting intact.\nconst replacements = [\n  [\"2024-03-19 Tuesday\", \"2024-03-20 Wednesday\"],\n  [\"913\u00f73=\", \"365\u00f78=\"],\n  [\"376\u00f72=\", \"496\u00f73=\"],\n  [\"557\u00f76=\", \"820\u00f74=\"],\n  [\"951\u00f75=\", \"184\u00f74=\"],\n  [\"941\u00f79=\", \"232\u00f76=\"],\n  [\"681\u00f78=\", \"168\u00f77=\"],\n  [\"231\u00f76=\", \"780\u00f79=\"],\n  [\"159\u00f78=\", \"653\u00f79=\"],\n  [\"211\u00f77=\", \"761\u00f76=\"],\n  [\"591\u00f77=\", \"397\u00f73=\"],\n  [\"829\u00f72=\", \"761\u00f74=\"],\n  [\"345\u00f76=\", \"321\u00f74=\"],\n  [\"930\u00f72=\", \"203\u00f76=\"],\n  [\"835\u00f77=\", \"104\u00f78=\"],\n  [\"849\u00f75=\", \"115\u00f73=\"],\n  [\"956\u00f72=\", \"889\u00f78=\"],\n  [\"170\u00f73=\", \"949\u00f73=\"],\n  [\"124\u00f73=\", \"124\u00f76=\"],\n  [\"209\u00f78=\", \"909\u00f73=\"],\n  [\"716\u00f72=\", \"234\u00f76=\"],\n  [\"814\u00f77=\", \"198\u00f76=\"],\n  [\"959\u00f72=\", \"456\u00f77=\"],\n  [\"365\u00f77=\", \"352\u00f79=\"],\n  [\"805\u00f79=\", \"417\u00f76=\"],\n  [\"510\u00f73=\", \"749\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division problems with their updated\n# values, as described by the diff. Every text run in the document is\n# mapped 1:1 (old text -> new text) and each old value is unique in the\n# document, so Find/Replace (wdReplaceAll) per pair is safe and keeps\n# each run's original formatting intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-19 Tuesday\", \"2024-03-20 Wednesday\"),\n    @(\"913\u00f73=\", \"365\u00f78=\"),\n    @(\"376\u00f72=\", \"496\u00f73=\"),\n    @(\"557\u00f76=\", \"820\u00f74=\"),\n    @(\"951\u00f75=\", \"184\u00f74=\"),\n    @(\"941\u00f79=\", \"232\u00f76=\"),\n    @(\"681\u00f78=\", \"168\u00f77=\"),\n    @(\"231\u00f76=\", \"780\u00f79=\"),\n    @(\"159\u00f78=\", \"653\u00f79=\"),\n    @(\"211\u00f77=\", \"761\u00f76=\"),\n    @(\"591\u00f77=\", \"397\u00f73=\"),\n    @(\"829\u00f72=\", \"761\u00f74=\"),\n    @(\"345\u00f76=\", \"321\u00f74=\"),\n    @(\"930\u00f72=\", \"203\u00f76=\"),\n    @(\"835\u00f77=\", \"104\u00f78=\"),\n    @(\"849\u00f75=\", \"115\u00f73=\"),\n    @(\"956\u00f72=\", \"889\u00f78=\"),\n    @(\"170\u00f73=\", \"949\u00f73=\"),\n    @(\"124\u00f73=\", \"124\u00f76=\"),\n    @(\"209\u00f78=\", \"909\u00f73=\"),\n    @(\"716\u00f72=\", \"234\u00f76=\"),\n    @(\"814\u00f77=\", \"198\u00f76=\"),\n    @(\"959\u00f72=\", \"456\u00f77=\"),\n    @(\"365\u00f77=\", \"352\u00f79=\"),\n    @(\"805\u00f79=\", \"417\u00f76=\"),\n    @(\"510\u00f73=\", \"749\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
